$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3316765.8
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 3316765.8
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 9950297.399999999
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -9950633.399999999
$ws.Range("H96").Value = 4118.8076
$ws.Range("I96").Value = 1894.2778
$ws.Range("J96").Value = 9124
$ws.Range("K96").Value = 5682.8334
$ws.Range("L96").Value = 27372
$ws.Range("M96").Value = -4309.8334
$ws.Range("N96").Value = -30118
$ws.Range("H100").Value = 1710.7
$ws.Range("I100").Value = 1352.1666
$ws.Range("J100").Value = 3144.8333
$ws.Range("K100").Value = 1352.1666
$ws.Range("L100").Value = 3144.8333
$ws.Range("M100").Value = -811.1666
$ws.Range("N100").Value = -4226.8333
$ws.Range("H103").Value = 850
$ws.Range("J103").Value = 1000
$ws.Range("L103").Value = 3000
$ws.Range("N103").Value = -4172
$ws.Range("H106").Value = 2665
$ws.Range("I106").Value = 1000
$ws.Range("K106").Value = 1000
$ws.Range("M106").Value = -369
$ws.Range("H107").Value = 2779.2
$ws.Range("I107").Value = 477.4
$ws.Range("J107").Value = 5081
$ws.Range("K107").Value = 477.4
$ws.Range("L107").Value = 5081
$ws.Range("M107").Value = 1442.6
$ws.Range("N107").Value = -8921
$ws.Range("H135").Value = 762.6
$ws.Range("I135").Value = 685.4545000000001
$ws.Range("K135").Value = 6169.0905
$ws.Range("M135").Value = -3634.0905
$ws.Range("H137").Value = 3673.4
$ws.Range("J137").Value = 4296.643
$ws.Range("L137").Value = 12889.929
$ws.Range("N137").Value = -17989.929
$ws.Range("H138").Value = 2879.6667
$ws.Range("J138").Value = 3579.5652
$ws.Range("L138").Value = 10738.6956
$ws.Range("N138").Value = -21018.6956
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3014
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H102").Value = 2166.4443
$ws.Range("I102").Value = 2298.375
$ws.Range("J102").Value = 1111
$ws.Range("K102").Value = 2298.375
$ws.Range("L102").Value = 1111
$ws.Range("M102").Value = -676.375
$ws.Range("N102").Value = -4355
$ws.Range("H122").Value = 1587.25
$ws.Range("I122").Value = 1471.1428
$ws.Range("K122").Value = 4413.428400000001
$ws.Range("M122").Value = -1963.428400000001
$ws.Range("H132").Value = 3336.9678
$ws.Range("I132").Value = 3336.9678
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10010.9034
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7480.903399999999
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 3014
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4347.2354
$ws.Range("J134").Value = 6517.2915
$ws.Range("L134").Value = 19551.8745
$ws.Range("N134").Value = -24621.8745
$ws.Range("H135").Value = 54997.5
$ws.Range("J135").Value = 54997.5
$ws.Range("L135").Value = 54997.5
$ws.Range("N135").Value = -65137.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2519.5918
$ws.Range("I31").Value = 1840.1923
$ws.Range("K31").Value = 1840.1923
$ws.Range("M31").Value = -1545.1923
$ws.Range("H34").Value = 2519.5918
$ws.Range("I34").Value = 1840.1923
$ws.Range("K34").Value = 1840.1923
$ws.Range("M34").Value = -1638.1923
$ws.Range("H51").Value = 9999.286
$ws.Range("J51").Value = 10832.5
$ws.Range("L51").Value = 10832.5
$ws.Range("N51").Value = -12304.5
$ws.Range("H59").Value = 20000
$ws.Range("J59").Value = 20000
$ws.Range("L59").Value = 20000
$ws.Range("N59").Value = -22290
$ws.Range("H60").Value = 8999.200000000001
$ws.Range("J60").Value = 8888.111000000001
$ws.Range("L60").Value = 8888.111000000001
$ws.Range("N60").Value = -9910.111000000001
$ws.Range("H61").Value = 9999.286
$ws.Range("J61").Value = 10832.5
$ws.Range("L61").Value = 10832.5
$ws.Range("N61").Value = -11528.5
$ws.Range("H99").Value = 20475.684
$ws.Range("I99").Value = 22317.076
$ws.Range("K99").Value = 22317.076
$ws.Range("M99").Value = -20819.076
$ws.Range("H126").Value = 20475.684
$ws.Range("I126").Value = 22317.076
$ws.Range("K126").Value = 66951.228
$ws.Range("M126").Value = -64481.228
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 615.63635
$ws.Range("I12").Value = 532.1667
$ws.Range("J12").Value = 715.8
$ws.Range("K12").Value = 1596.5001
$ws.Range("L12").Value = 2147.4
$ws.Range("M12").Value = -1423.5001
$ws.Range("N12").Value = -2493.4
$ws.Range("H92").Value = 989.5
$ws.Range("I92").Value = 934.1429000000001
$ws.Range("J92").Value = 1044.8572
$ws.Range("K92").Value = 2802.4287
$ws.Range("L92").Value = 3134.5716
$ws.Range("M92").Value = -1554.4287
$ws.Range("N92").Value = -5630.571599999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2673.8333
$ws.Range("I102").Value = 2673.8333
$ws.Range("K102").Value = 2673.8333
$ws.Range("M102").Value = -1051.8333
$ws.Range("H132").Value = 4142.9546
$ws.Range("I132").Value = 4396.8237
$ws.Range("J132").Value = 3279.8
$ws.Range("K132").Value = 13190.4711
$ws.Range("L132").Value = 9839.400000000001
$ws.Range("M132").Value = -10660.4711
$ws.Range("N132").Value = -14899.4
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1212.1333
$ws.Range("I16").Value = 1314
$ws.Range("J16").Value = 550
$ws.Range("K16").Value = 1314
$ws.Range("L16").Value = 550
$ws.Range("M16").Value = -1144
$ws.Range("N16").Value = -890
$ws.Range("H22").Value = 1622.2222
$ws.Range("I22").Value = 1120
$ws.Range("K22").Value = 1120
$ws.Range("M22").Value = -825
$ws.Range("H27").Value = 1622.2222
$ws.Range("I27").Value = 1120
$ws.Range("K27").Value = 1120
$ws.Range("M27").Value = -1013
$ws.Range("H55").Value = 323.73077
$ws.Range("I55").Value = 259.77777
$ws.Range("K55").Value = 259.77777
$ws.Range("M55").Value = -86.77776999999998
$ws.Range("H61").Value = 83421170
$ws.Range("I61").Value = 166671170
$ws.Range("J61").Value = 171164
$ws.Range("K61").Value = 166671170
$ws.Range("L61").Value = 171164
$ws.Range("M61").Value = -166670968
$ws.Range("N61").Value = -171568
$ws.Range("H93").Value = 19799.666
$ws.Range("J93").Value = 84583
$ws.Range("L93").Value = 84583
$ws.Range("N93").Value = -87079
$ws.Range("H113").Value = 83421170
$ws.Range("I113").Value = 166671170
$ws.Range("J113").Value = 171164
$ws.Range("K113").Value = 166671170
$ws.Range("L113").Value = 171164
$ws.Range("M113").Value = -166669000
$ws.Range("N113").Value = -175504
$ws.Range("H132").Value = 4856.2856
$ws.Range("I132").Value = 3666.3333
$ws.Range("J132").Value = 5748.75
$ws.Range("K132").Value = 10998.9999
$ws.Range("L132").Value = 17246.25
$ws.Range("M132").Value = -8468.999899999999
$ws.Range("N132").Value = -22306.25
$ws.Range("H136").Value = 3999
$ws.Range("I136").Value = 3999
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 11997
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -9447
$ws.Range("N136").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 28639.182
$ws.Range("I45").Value = 15783.5
$ws.Range("J45").Value = 31496
$ws.Range("K45").Value = 15783.5
$ws.Range("L45").Value = 31496
$ws.Range("M45").Value = -15292.5
$ws.Range("N45").Value = -32478
$ws.Range("H122").Value = 1817.6562
$ws.Range("I122").Value = 1316.8823
$ws.Range("K122").Value = 3950.6469
$ws.Range("M122").Value = -1500.6469
$ws.Range("H132").Value = 7420.037
$ws.Range("I132").Value = 8159.2383
$ws.Range("J132").Value = 4832.8335
$ws.Range("K132").Value = 24477.7149
$ws.Range("L132").Value = 14498.5005
$ws.Range("M132").Value = -21947.7149
$ws.Range("N132").Value = -19558.5005
$ws.Range("H136").Value = 6342.839
$ws.Range("I136").Value = 4178.4443
$ws.Range("J136").Value = 20952.5
$ws.Range("K136").Value = 12535.3329
$ws.Range("L136").Value = 62857.5
$ws.Range("M136").Value = -9985.332900000001
$ws.Range("N136").Value = -67957.5
